$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8163.6665
$ws.Range("I51").Value = 6994.6665
$ws.Range("K51").Value = 6994.6665
$ws.Range("M51").Value = -6510.6665
$ws.Range("H74").Value = 10932.467
$ws.Range("I74").Value = 10192.077
$ws.Range("K74").Value = 10192.077
$ws.Range("M74").Value = -9256.076999999999
$ws.Range("H77").Value = 10932.467
$ws.Range("I77").Value = 10192.077
$ws.Range("K77").Value = 50960.38499999999
$ws.Range("M77").Value = -46280.38499999999
$ws.Range("H125").Value = 8565.857
$ws.Range("I125").Value = 7562.2
$ws.Range("J125").Value = 11075
$ws.Range("K125").Value = 68059.8
$ws.Range("L125").Value = 99675
$ws.Range("M125").Value = -65599.8
$ws.Range("N125").Value = -104595
$ws.Range("H132").Value = 6120.8
$ws.Range("I132").Value = 7136.5884
$ws.Range("J132").Value = 3962.25
$ws.Range("K132").Value = 21409.7652
$ws.Range("L132").Value = 11886.75
$ws.Range("M132").Value = -18879.7652
$ws.Range("N132").Value = -16946.75
$ws.Range("H137").Value = 4860.7144
$ws.Range("I137").Value = 5068
$ws.Range("J137").Value = 4804.1816
$ws.Range("K137").Value = 15204
$ws.Range("L137").Value = 14412.5448
$ws.Range("M137").Value = -12654
$ws.Range("N137").Value = -19512.5448
$ws.Range("H138").Value = 4409.1274
$ws.Range("J138").Value = 4582.638
$ws.Range("L138").Value = 13747.914
$ws.Range("N138").Value = -24027.914

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1254999.2
$ws.Range("I110").Value = 1254999.2
$ws.Range("K110").Value = 1254999.2
$ws.Range("M110").Value = -1252954.2
$ws.Range("H122").Value = 5715.9443
$ws.Range("I122").Value = 5420.9287
$ws.Range("K122").Value = 16262.7861
$ws.Range("M122").Value = -13812.7861
$ws.Range("H132").Value = 4614.5
$ws.Range("I132").Value = 1241.5
$ws.Range("K132").Value = 3724.5
$ws.Range("M132").Value = -1194.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 299184
$ws.Range("J42").Value = 299184
$ws.Range("L42").Value = 299184
$ws.Range("N42").Value = -299840
$ws.Range("H105").Value = 68690.92999999999
$ws.Range("I105").Value = 250972.5
$ws.Range("J105").Value = 2406.7273
$ws.Range("K105").Value = 250972.5
$ws.Range("L105").Value = 2406.7273
$ws.Range("M105").Value = -249225.5
$ws.Range("N105").Value = -5900.7273
$ws.Range("H123").Value = 74622.5
$ws.Range("J123").Value = 74622.5
$ws.Range("L123").Value = 74622.5
$ws.Range("N123").Value = -84422.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 273919.03
$ws.Range("I58").Value = 419339.47
$ws.Range("J58").Value = 5450.5386
$ws.Range("K58").Value = 419339.47
$ws.Range("L58").Value = 5450.5386
$ws.Range("M58").Value = -419136.47
$ws.Range("N58").Value = -5856.5386
$ws.Range("H96").Value = 39906
$ws.Range("J96").Value = 39906
$ws.Range("L96").Value = 39906
$ws.Range("N96").Value = -45398
$ws.Range("H99").Value = 6109.8184
$ws.Range("I99").Value = 5052.25
$ws.Range("K99").Value = 5052.25
$ws.Range("M99").Value = -3554.25
$ws.Range("H122").Value = 3858.4443
$ws.Range("I122").Value = 2178
$ws.Range("K122").Value = 6534
$ws.Range("M122").Value = -4084
$ws.Range("H126").Value = 6109.8184
$ws.Range("I126").Value = 5052.25
$ws.Range("K126").Value = 15156.75
$ws.Range("M126").Value = -12686.75
$ws.Range("H136").Value = 273919.03
$ws.Range("I136").Value = 419339.47
$ws.Range("J136").Value = 5450.5386
$ws.Range("K136").Value = 1258018.41
$ws.Range("L136").Value = 16351.6158
$ws.Range("M136").Value = -1255468.41
$ws.Range("N136").Value = -21451.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 16.913044
$ws.Range("I2").Value = 12.111111
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 72.66666599999999
$ws.Range("L2").Value = 120
$ws.Range("M2").Value = 40.33333400000001
$ws.Range("N2").Value = -346
$ws.Range("H9").Value = 2625118.5
$ws.Range("J9").Value = 3000135.5
$ws.Range("L9").Value = 9000406.5
$ws.Range("N9").Value = -9000854.5
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 200
$ws.Range("K20").Value = 600
$ws.Range("M20").Value = -373
$ws.Range("H21").Value = 209.5
$ws.Range("I21").Value = 168
$ws.Range("J21").Value = 251
$ws.Range("K21").Value = 504
$ws.Range("L21").Value = 753
$ws.Range("M21").Value = -331
$ws.Range("N21").Value = -1099
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H23").Value = 264.3
$ws.Range("J23").Value = 262.33334
$ws.Range("L23").Value = 787.0000200000001
$ws.Range("N23").Value = -1257.00002
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H57").Value = 12000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H130").Value = 4565.5
$ws.Range("I130").Value = 1623.75
$ws.Range("J130").Value = 10449
$ws.Range("K130").Value = 4871.25
$ws.Range("L130").Value = 31347
$ws.Range("M130").Value = 148.75
$ws.Range("N130").Value = -41387
$ws.Range("H131").Value = 3162.8333
$ws.Range("I131").Value = 1362.4286
$ws.Range("J131").Value = 4308.5454
$ws.Range("K131").Value = 4087.2858
$ws.Range("L131").Value = 12925.6362
$ws.Range("M131").Value = 952.7142000000003
$ws.Range("N131").Value = -23005.6362
$ws.Range("H136").Value = 5919
$ws.Range("I136").Value = 4775.5
$ws.Range("J136").Value = 9349.5
$ws.Range("K136").Value = 14326.5
$ws.Range("L136").Value = 28048.5
$ws.Range("M136").Value = -9226.5
$ws.Range("N136").Value = -38248.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 532.8
$ws.Range("I22").Value = 480
$ws.Range("K22").Value = 480
$ws.Range("M22").Value = -185
$ws.Range("H27").Value = 532.8
$ws.Range("I27").Value = 480
$ws.Range("K27").Value = 480
$ws.Range("M27").Value = -373
$ws.Range("H100").Value = 179800.17
$ws.Range("I100").Value = 201760.2
$ws.Range("J100").Value = 70000
$ws.Range("K100").Value = 201760.2
$ws.Range("L100").Value = 70000
$ws.Range("M100").Value = -201219.2
$ws.Range("N100").Value = -71082
$ws.Range("H106").Value = 14185
$ws.Range("J106").Value = 14185
$ws.Range("L106").Value = 14185
$ws.Range("N106").Value = -16709
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 5560.2666
$ws.Range("I132").Value = 4601.3335
$ws.Range("J132").Value = 5800
$ws.Range("K132").Value = 13804.0005
$ws.Range("L132").Value = 17400
$ws.Range("M132").Value = -11274.0005
$ws.Range("N132").Value = -22460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9600.857
$ws.Range("J41").Value = 9271.333000000001
$ws.Range("L41").Value = 9271.333000000001
$ws.Range("N41").Value = -10051.333
$ws.Range("H81").Value = 20947.77
$ws.Range("J81").Value = 4824
$ws.Range("L81").Value = 9648
$ws.Range("N81").Value = -11770
$ws.Range("H84").Value = 20947.77
$ws.Range("J84").Value = 4824
$ws.Range("L84").Value = 48240
$ws.Range("N84").Value = -58848
$ws.Range("H96").Value = 112730.78
$ws.Range("I96").Value = 201835.6
$ws.Range("J96").Value = 1349.75
$ws.Range("K96").Value = 201835.6
$ws.Range("L96").Value = 1349.75
$ws.Range("M96").Value = -200462.6
$ws.Range("N96").Value = -4095.75
$ws.Range("H100").Value = 356.7857
$ws.Range("I100").Value = 252.27272
$ws.Range("K100").Value = 504.54544
$ws.Range("M100").Value = 36.45456000000001
$ws.Range("H121").Value = 61806.668
$ws.Range("J121").Value = 61806.668
$ws.Range("L121").Value = 61806.668
$ws.Range("N121").Value = -65300.668
